$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, date range) ---
$ws.Range("A8").Value = "Volume 32   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/12/2025  Through  5/18/2025"

# --- Simple numeric cell updates (rows 15-19, 21, 24-28; no type change) ---
$ws.Range("L15").Value = -40
$ws.Range("N15").Value = -62.5

$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 36.363636363636
$ws.Range("I16").Value = 67
$ws.Range("J16").Value = 70
$ws.Range("K16").Value = -4.285714285714
$ws.Range("L16").Value = -4.285714285714
$ws.Range("M16").Value = 34
$ws.Range("N16").Value = -87.040618955512

$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = -41.379310344827
$ws.Range("I17").Value = 69
$ws.Range("J17").Value = 69
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = -5.479452054794
$ws.Range("M17").Value = 91.666666666666
$ws.Range("N17").Value = -35.514018691588

$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = 14.285714285714
$ws.Range("F18").Value = 27
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = 22.727272727272
$ws.Range("I18").Value = 111
$ws.Range("J18").Value = 95
$ws.Range("K18").Value = 16.842105263157
$ws.Range("L18").Value = -1.769911504424
$ws.Range("M18").Value = 12.121212121212
$ws.Range("N18").Value = -89.16015625

$ws.Range("C19").Value = 26
$ws.Range("D19").Value = 37
$ws.Range("E19").Value = -29.729729729729
$ws.Range("F19").Value = 126
$ws.Range("G19").Value = 116
$ws.Range("H19").Value = 8.620689655172
$ws.Range("I19").Value = 569
$ws.Range("J19").Value = 576
$ws.Range("K19").Value = -1.215277777777
$ws.Range("L19").Value = -4.849498327759
$ws.Range("M19").Value = 23.965141612200
$ws.Range("N19").Value = -59.559346126510

$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 52
$ws.Range("E21").Value = -25
$ws.Range("F21").Value = 188
$ws.Range("G21").Value = 180
$ws.Range("H21").Value = 4.444444444444
$ws.Range("I21").Value = 837
$ws.Range("J21").Value = 837
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = -8.823529411764
$ws.Range("M21").Value = 24.553571428571
$ws.Range("N21").Value = -80.959963603275

$ws.Range("C24").Value = 58
$ws.Range("D24").Value = 58
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 177
$ws.Range("G24").Value = 208
$ws.Range("H24").Value = -14.903846153846
$ws.Range("I24").Value = 1171
$ws.Range("J24").Value = 1078
$ws.Range("K24").Value = 8.627087198515
$ws.Range("L24").Value = 1.473136915077
$ws.Range("M24").Value = 96.147403685092

$ws.Range("C25").Value = 41
$ws.Range("D25").Value = 52
$ws.Range("E25").Value = -21.153846153846
$ws.Range("F25").Value = 147
$ws.Range("G25").Value = 184
$ws.Range("H25").Value = -20.108695652173
$ws.Range("I25").Value = 995
$ws.Range("J25").Value = 952
$ws.Range("K25").Value = 4.516806722689
$ws.Range("L25").Value = -0.698602794411

$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 32
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = 14.285714285714
$ws.Range("I26").Value = 119
$ws.Range("J26").Value = 129
$ws.Range("K26").Value = -7.751937984496
$ws.Range("L26").Value = -11.851851851851
$ws.Range("M26").Value = -10.526315789473

$ws.Range("L27").Value = -33.333333333333

$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 4
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = 12
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 32
$ws.Range("J28").Value = 52
$ws.Range("K28").Value = -38.461538461538
$ws.Range("L28").Value = -17.948717948717

# --- Row 20: C/D/E flip from numbers to text ("0"/"0"/"***.*"); F/G/H/L/N numeric updates ---
$ws.Range("C20").Value = "'0"
$ws.Range("A20").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Range("D20").Value = "'0"
$ws.Range("A20").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").Value = "'***.*"
$ws.Range("A20").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 50
$ws.Range("L20").Value = -68.421052631578
$ws.Range("N20").Value = -98.646616541353

# --- Row 22: D/E flip from text to numbers; G/H/J/K numeric updates ---
$ws.Range("D22").Value = 2
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -75
$ws.Range("J22").Value = 16
$ws.Range("K22").Value = -25

# --- Row 23: D/E flip from numbers to text ("0"/"***.*"); F/H/M numeric updates ---
$ws.Range("D23").Value = "'0"
$ws.Range("C23").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Value = "'***.*"
$ws.Range("N23").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null
$ws.Range("F23").Value = 1
$ws.Range("H23").Value = 0
$ws.Range("M23").Value = 8.333333333333

# --- Row 31: D/E flip from text to numbers; G/H/J/K/L numeric updates ---
$ws.Range("D31").Value = 1
$ws.Range("D31").NumberFormat = "#,##0"
$ws.Range("E31").Value = -100
$ws.Range("E31").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = -50
$ws.Range("J31").Value = 13
$ws.Range("K31").Value = -38.461538461538
$ws.Range("L31").Value = 60

Write-Host "Edit complete"